$d = $word.ActiveDocument

# Update the date line (unique text, safe to use Find/Replace)
$d.Content.Find.Execute("2023-08-11 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-12 Saturday", 2)

# Update the division problems table. Cell text is set directly (by row/col)
# rather than via global Find/Replace, because several of the new values
# coincide with pre-existing / other target values elsewhere in the table
# (e.g. "37÷2=" -> "79÷9=" while the original "79÷9=" -> "21÷6="); a
# document-wide replace-all could cross-match those and corrupt entries.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "79÷9="
$t.Cell(1, 2).Range.Text = "50÷2="
$t.Cell(1, 3).Range.Text = "65÷5="
$t.Cell(1, 4).Range.Text = "15÷5="
$t.Cell(1, 5).Range.Text = "54÷3="

$t.Cell(5, 1).Range.Text = "72÷5="
$t.Cell(5, 2).Range.Text = "13÷7="
$t.Cell(5, 3).Range.Text = "42÷2="
$t.Cell(5, 4).Range.Text = "16÷5="
$t.Cell(5, 5).Range.Text = "18÷8="

$t.Cell(9, 1).Range.Text = "22÷6="
$t.Cell(9, 2).Range.Text = "71÷6="
$t.Cell(9, 3).Range.Text = "21÷6="
$t.Cell(9, 4).Range.Text = "89÷6="
$t.Cell(9, 5).Range.Text = "25÷6="

$t.Cell(13, 1).Range.Text = "12÷6="
$t.Cell(13, 2).Range.Text = "57÷9="
$t.Cell(13, 3).Range.Text = "95÷5="
$t.Cell(13, 4).Range.Text = "82÷4="
$t.Cell(13, 5).Range.Text = "46÷9="

$t.Cell(17, 1).Range.Text = "49÷5="
$t.Cell(17, 2).Range.Text = "32÷6="
$t.Cell(17, 3).Range.Text = "34÷2="
$t.Cell(17, 4).Range.Text = "49÷4="
$t.Cell(17, 5).Range.Text = "48÷7="
